$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns F (JSHIR), I (Telegram), J (Telefon) and K (Sana) contain values that
# look like numbers/dates to Excel's automatic type detection (long digit
# strings, a "+" prefixed phone number, and an ISO date string). Force those
# columns to Text format first so the values are stored verbatim instead of
# being reinterpreted as numbers or date serials.
$ws.Range("F20:F21").NumberFormat = "@"
$ws.Range("I20:J21").NumberFormat = "@"
$ws.Range("K20:K21").NumberFormat = "@"

$rows = @(
    @{ Row = 20; Values = @(
        "Xudoyorov Muslimjon Mominjon ogli",
        "Yurisprudensiya",
        "O'zbek tili",
        "Sirtqi",
        "AD7761080",
        "31103914340034",
        "Fargona viloyati",
        "Marg" + [char]0x02BB + "ilon tumani",
        "998901669999",
        "+998916588000",
        "2025-04-27"
    ) },
    @{ Row = 21; Values = @(
        "Adizov Ismoiljon",
        "Yurisprudensiya",
        "O'zbek tili",
        "Sirtqi",
        "AD4325461",
        "32804881070096",
        "Toshkent shahri",
        "Mirzo Ulug" + [char]0x02BB + "bek tumani",
        "998903490733",
        "+998936578677",
        "2025-04-27"
    ) }
)

foreach ($rowInfo in $rows) {
    $r = $rowInfo.Row
    $colIndex = 1
    foreach ($val in $rowInfo.Values) {
        $ws.Cells.Item($r, $colIndex).Value = $val
        $colIndex++
    }
}
